# 13.1.3.xlsx update:
#  - Replace the old "484 / 1.5.4 ..." text in A1 with the new 13.1.3 header text
#  - The shared-strings table is compacted automatically on save, so simply
#    writing the desired final text/values into cells is enough; unused
#    strings (old "484" text, old "1.5.4 ..." text) drop out on their own
#    once nothing references them any more.
#  - Extend the data table with columns 2020-2023 (E:H), mirroring the
#    formatting of the existing 2019 column (D) via Range.Copy, then
#    overwriting the copied values with the real figures.
#  - Turn the D4 "484" text cell into a real number 484.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: new 13.1.3 wording -----------------------------------
$ws.Range("A1").Value = "13.1.3 Кырсыктардын кооптуулугун азайтуунун улуттук стратегияларына ылайык, кырсыктардын кооптуулугун азайтуунун жергиликтүү стратегияларын кабыл алган жана ишке ашырган жергиликтүү бийлик органдарынын үлүшү"

# --- New year columns: headers 2020..2023 in E3:H3 ---------------------
$ws.Range("D3").Copy($ws.Range("E3:H3"))
$ws.Range("E3").Value = 2020
$ws.Range("F3").Value = 2021
$ws.Range("G3").Value = 2022
$ws.Range("H3").Value = 2023

# --- Row 4: "number of local governments" = 484, now a real number -----
$ws.Range("D4").Value = 484
$ws.Range("D4").Copy($ws.Range("E4:H4"))
$ws.Range("E4").Value = 484
$ws.Range("F4").Value = 484
$ws.Range("G4").Value = 484
$ws.Range("H4").Value = 484

# --- Row 5: proportion (%) ---------------------------------------------
$ws.Range("D5").Copy($ws.Range("E5:H5"))
$ws.Range("E5").Value = 13.2
$ws.Range("F5").Value = 21.5
$ws.Range("G5").Value = 34.5
$ws.Range("H5").Value = 40.53

# --- Row 6: number of local governments that adopted DRR strategies ----
$ws.Range("D6").Copy($ws.Range("E6:H6"))
$ws.Range("E6").Value = 67
$ws.Range("F6").Value = 104
$ws.Range("G6").Value = 167
$ws.Range("H6").Value = 169
